$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.721.89"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "2.287.59"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "115.69"
$ws.Range("E5").Value = "  +11.72%  "
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.72"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.37"
$ws.Range("E10").Value = "  +8.05%  "
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0946"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.93"
$ws.Range("E12").Value = "  +11.50%  "
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.92"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "2.631.49"
$ws.Range("E15").Value = "  -0.29%  "

$ws.Range("E16").Value = "  +2.96%  "

$ws.Range("D17").Value = "2.287.38"
$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "43.619.94"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = "  +12.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.34"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("E23").Value = "  +8.93%  "
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.64"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.66"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "42.12"
$ws.Range("E28").Value = "  +5.55%  "
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").Value = "3.39"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "172.85"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0931"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "21.62"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  +4.84%  "
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.66"
$ws.Range("E36").Value = "  -5.64%  "
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0358"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.78"
$ws.Range("E39").Value = "  +6.24%  "
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "14.67"
$ws.Range("E40").Value = "  +19.56%  "
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "2.44"
$ws.Range("E41").Value = "  +4.53%  "
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "74.87"
$ws.Range("E42").Value = "  +14.12%  "
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.241"
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "6.38"
$ws.Range("E44").Value = "  +21.57%  "
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "1.39"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "8.70"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "1.27"
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "103.04"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.459"
$ws.Range("E51").Value = "  +1.57%  "
$ws.Range("D51").Style = "Normal"
